$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.535856000000001
$ws.Range("H2").Value = 25.607568
$ws.Range("I2").Value = 0.36987004643386
$ws.Range("J2").Value = 0.36987004643386
$ws.Range("M2").Value = 2.207767333333333
$ws.Range("N2").Value = 6.623302000000001
$ws.Range("O2").Value = 0.07574879669493666
$ws.Range("P2").Value = 0.07574879669493666
$ws.Range("Q2").Value = 18.84518403883733
$ws.Range("R2").Value = 169.606656349536
$ws.Range("S2").Value = 0.02801721095086524
$ws.Range("T2").Value = 0.02801721095086524
$ws.Range("G3").Value = 8.535856000000001
$ws.Range("H3").Value = 25.607568
$ws.Range("I3").Value = 0.36987004643386
$ws.Range("J3").Value = 0.36987004643386
$ws.Range("O3").Value = 0.6498413459298955
$ws.Range("P3").Value = 0.6498413459298955
$ws.Range("Q3").Value = 161.67094784904
$ws.Range("R3").Value = 1455.03853064136
$ws.Range("S3").Value = 0.2403568487937325
$ws.Range("T3").Value = 0.2403568487937325
$ws.Range("G4").Value = 8.535856000000001
$ws.Range("H4").Value = 25.607568
$ws.Range("I4").Value = 0.36987004643386
$ws.Range("J4").Value = 0.36987004643386
$ws.Range("M4").Value = 7.997924
$ws.Range("N4").Value = 23.993772
$ws.Range("O4").Value = 0.2744098573751678
$ws.Range("P4").Value = 0.2744098573751678
$ws.Range("Q4").Value = 68.26912756294401
$ws.Range("R4").Value = 614.422148066496
$ws.Range("S4").Value = 0.1014959866892622
$ws.Range("T4").Value = 0.1014959866892622
$ws.Range("I5").Value = 0.2236685002562326
$ws.Range("J5").Value = 0.2236685002562326
$ws.Range("M5").Value = 2.207767333333333
$ws.Range("N5").Value = 6.623302000000001
$ws.Range("O5").Value = 0.07574879669493666
$ws.Range("P5").Value = 0.07574879669493666
$ws.Range("Q5").Value = 11.39609463285689
$ws.Range("R5").Value = 102.564851695712
$ws.Range("S5").Value = 0.01694261975297075
$ws.Range("T5").Value = 0.01694261975297075
$ws.Range("I6").Value = 0.2236685002562326
$ws.Range("J6").Value = 0.2236685002562326
$ws.Range("O6").Value = 0.6498413459298955
$ws.Range("P6").Value = 0.6498413459298955
$ws.Range("S6").Value = 0.1453490392486314
$ws.Range("T6").Value = 0.1453490392486314
$ws.Range("I7").Value = 0.2236685002562326
$ws.Range("J7").Value = 0.2236685002562326
$ws.Range("M7").Value = 7.997924
$ws.Range("N7").Value = 23.993772
$ws.Range("O7").Value = 0.2744098573751678
$ws.Range("P7").Value = 0.2744098573751678
$ws.Range("Q7").Value = 41.28383339778133
$ws.Range("R7").Value = 371.554500580032
$ws.Range("S7").Value = 0.06137684125463048
$ws.Range("T7").Value = 0.06137684125463048
$ws.Range("G8").Value = 5.335438
$ws.Range("H8").Value = 16.006314
$ws.Range("I8").Value = 0.2311916579666972
$ws.Range("J8").Value = 0.2311916579666973
$ws.Range("M8").Value = 2.207767333333333
$ws.Range("N8").Value = 6.623302000000001
$ws.Range("O8").Value = 0.07574879669493666
$ws.Range("P8").Value = 0.07574879669493666
$ws.Range("Q8").Value = 11.77940572542533
$ws.Range("R8").Value = 106.014651528828
$ws.Range("S8").Value = 0.01751248989688468
$ws.Range("T8").Value = 0.01751248989688468
$ws.Range("G9").Value = 5.335438
$ws.Range("H9").Value = 16.006314
$ws.Range("I9").Value = 0.2311916579666972
$ws.Range("J9").Value = 0.2311916579666973
$ws.Range("O9").Value = 0.6498413459298955
$ws.Range("P9").Value = 0.6498413459298955
$ws.Range("Q9").Value = 101.05434283917
$ws.Range("R9").Value = 909.48908555253
$ws.Range("S9").Value = 0.1502378981808426
$ws.Range("T9").Value = 0.1502378981808426
$ws.Range("G10").Value = 5.335438
$ws.Range("H10").Value = 16.006314
$ws.Range("I10").Value = 0.2311916579666972
$ws.Range("J10").Value = 0.2311916579666973
$ws.Range("M10").Value = 7.997924
$ws.Range("N10").Value = 23.993772
$ws.Range("O10").Value = 0.2744098573751678
$ws.Range("P10").Value = 0.2744098573751678
$ws.Range("Q10").Value = 42.672427630712
$ws.Range("R10").Value = 384.051848676408
$ws.Range("S10").Value = 0.06344126988896998
$ws.Range("T10").Value = 0.06344126988896998
$ws.Range("G11").Value = 4.044874
$ws.Range("H11").Value = 12.134622
$ws.Range("I11").Value = 0.1752697953432102
$ws.Range("J11").Value = 0.1752697953432102
$ws.Range("M11").Value = 2.207767333333333
$ws.Range("N11").Value = 6.623302000000001
$ws.Range("O11").Value = 0.07574879669493666
$ws.Range("P11").Value = 0.07574879669493666
$ws.Range("Q11").Value = 8.930140684649334
$ws.Range("R11").Value = 80.37126616184401
$ws.Range("S11").Value = 0.01327647609421598
$ws.Range("T11").Value = 0.01327647609421598
$ws.Range("G12").Value = 4.044874
$ws.Range("H12").Value = 12.134622
$ws.Range("I12").Value = 0.1752697953432102
$ws.Range("J12").Value = 0.1752697953432102
$ws.Range("O12").Value = 0.6498413459298955
$ws.Range("P12").Value = 0.6498413459298955
$ws.Range("Q12").Value = 76.61078320790999
$ws.Range("R12").Value = 689.49704887119
$ws.Range("S12").Value = 0.113897559706689
$ws.Range("T12").Value = 0.113897559706689
$ws.Range("G13").Value = 4.044874
$ws.Range("H13").Value = 12.134622
$ws.Range("I13").Value = 0.1752697953432102
$ws.Range("J13").Value = 0.1752697953432102
$ws.Range("M13").Value = 7.997924
$ws.Range("N13").Value = 23.993772
$ws.Range("O13").Value = 0.2744098573751678
$ws.Range("P13").Value = 0.2744098573751678
$ws.Range("Q13").Value = 32.350594841576
$ws.Range("R13").Value = 291.155353574184
$ws.Range("S13").Value = 0.04809575954230516
$ws.Range("T13").Value = 0.04809575954230516
